$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 266450.75
$ws.Range("I4").Value = 353600.34
$ws.Range("J4").Value = 5002
$ws.Range("K4").Value = 353600.34
$ws.Range("L4").Value = 5002
$ws.Range("M4").Value = -353486.34
$ws.Range("N4").Value = -5230
# Row 70
$ws.Range("H70").Value = 2957.5925
$ws.Range("J70").Value = 3194.125
$ws.Range("L70").Value = 9582.375
$ws.Range("N70").Value = -10122.375
# Row 73
$ws.Range("H73").Value = 2957.5925
$ws.Range("J73").Value = 3194.125
$ws.Range("L73").Value = 9582.375
$ws.Range("N73").Value = -11454.375
# Row 107
$ws.Range("H107").Value = 907397.9
$ws.Range("J107").Value = 9500
$ws.Range("L107").Value = 9500
$ws.Range("N107").Value = -13340
# Row 137
$ws.Range("H137").Value = 1886.9584
$ws.Range("J137").Value = 1733.3334
$ws.Range("L137").Value = 5200.0002
$ws.Range("N137").Value = -10300.0002
# Row 138
$ws.Range("H138").Value = 3073.3333
$ws.Range("I138").Value = 1683.9231
$ws.Range("J138").Value = 4363.5
$ws.Range("K138").Value = 5051.7693
$ws.Range("L138").Value = 13090.5
$ws.Range("M138").Value = 88.23070000000007
$ws.Range("N138").Value = -23370.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3252.9648
$ws.Range("I32").Value = 3458.9546
$ws.Range("K32").Value = 3458.9546
$ws.Range("M32").Value = -3171.9546
# Row 61
$ws.Range("H61").Value = 3898.1724
$ws.Range("I61").Value = 2954.75
$ws.Range("J61").Value = 5059.3076
$ws.Range("K61").Value = 2954.75
$ws.Range("L61").Value = 5059.3076
$ws.Range("M61").Value = -2742.75
$ws.Range("N61").Value = -5483.3076
# Row 74
$ws.Range("H74").Value = 1642.48
$ws.Range("I74").Value = 1450.6842
$ws.Range("K74").Value = 1450.6842
$ws.Range("M74").Value = -576.6841999999999
# Row 77
$ws.Range("H77").Value = 1642.48
$ws.Range("I77").Value = 1450.6842
$ws.Range("K77").Value = 7253.420999999999
$ws.Range("M77").Value = -2885.420999999999
# Row 135
$ws.Range("H135").Value = 77000
$ws.Range("J135").Value = 77000
$ws.Range("L135").Value = 77000
$ws.Range("N135").Value = -87140
# Row 136
$ws.Range("H136").Value = 3898.1724
$ws.Range("I136").Value = 2954.75
$ws.Range("J136").Value = 5059.3076
$ws.Range("K136").Value = 8864.25
$ws.Range("L136").Value = 15177.9228
$ws.Range("M136").Value = -6314.25
$ws.Range("N136").Value = -20277.9228
# Row 139
$ws.Range("H139").Value = 87944.22
$ws.Range("J139").Value = 87944.22
$ws.Range("L139").Value = 87944.22
$ws.Range("N139").Value = -98224.22

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3757
$ws.Range("I58").Value = 3908.6
$ws.Range("J58").Value = 2999
$ws.Range("K58").Value = 3908.6
$ws.Range("L58").Value = 2999
$ws.Range("M58").Value = -3705.6
$ws.Range("N58").Value = -3405
# Row 132
$ws.Range("H132").Value = 4170.615
$ws.Range("I132").Value = 4185.5835
$ws.Range("K132").Value = 12556.7505
$ws.Range("M132").Value = -10026.7505
# Row 136
$ws.Range("H136").Value = 3757
$ws.Range("I136").Value = 3908.6
$ws.Range("J136").Value = 2999
$ws.Range("K136").Value = 11725.8
$ws.Range("L136").Value = 8997
$ws.Range("M136").Value = -9175.799999999999
$ws.Range("N136").Value = -14097

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 2833.1052
$ws.Range("I14").Value = 2833.1052
$ws.Range("K14").Value = 8499.3156
$ws.Range("M14").Value = -8326.3156
# Row 41
$ws.Range("H41").Value = 484
$ws.Range("I41").Value = 210
$ws.Range("J41").Value = 666.6667
$ws.Range("K41").Value = 630
$ws.Range("L41").Value = 2000.0001
$ws.Range("M41").Value = -292
$ws.Range("N41").Value = -2676.0001
# Row 107
$ws.Range("H107").Value = 1657.2264
$ws.Range("J107").Value = 1797.3405
$ws.Range("L107").Value = 5392.0215
$ws.Range("N107").Value = -9232.021499999999
# Row 131
$ws.Range("H131").Value = 1857.0526
$ws.Range("I131").Value = 1236.6666
$ws.Range("K131").Value = 3709.9998
$ws.Range("M131").Value = 1330.0002
# Row 137
$ws.Range("H137").Value = 1898.7191
$ws.Range("J137").Value = 1971.7949
$ws.Range("L137").Value = 5915.384700000001
$ws.Range("N137").Value = -16115.3847

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2294.4285
$ws.Range("I80").Value = 2011
$ws.Range("K80").Value = 2011
$ws.Range("M80").Value = -1013
# Row 83
$ws.Range("H83").Value = 2294.4285
$ws.Range("I83").Value = 2011
$ws.Range("K83").Value = 10055
$ws.Range("M83").Value = -5063
# Row 126
$ws.Range("H126").Value = 3630.6667
$ws.Range("I126").Value = 4170.2
$ws.Range("J126").Value = 2956.25
$ws.Range("K126").Value = 12510.6
$ws.Range("L126").Value = 8868.75
$ws.Range("M126").Value = -10040.6
$ws.Range("N126").Value = -13808.75
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 17
$ws.Range("H17").Value = 3387.6924
$ws.Range("I17").Value = 2677.5
$ws.Range("J17").Value = 3703.3333
$ws.Range("K17").Value = 2677.5
$ws.Range("L17").Value = 3703.3333
$ws.Range("M17").Value = -2507.5
$ws.Range("N17").Value = -4043.3333
# Row 46
$ws.Range("H46").Value = 2628.25
$ws.Range("I46").Value = 1081.75
$ws.Range("J46").Value = 4174.75
$ws.Range("K46").Value = 1081.75
$ws.Range("L46").Value = 4174.75
$ws.Range("M46").Value = -893.75
$ws.Range("N46").Value = -4550.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 662.4706
$ws.Range("I100").Value = 577.93335
$ws.Range("K100").Value = 1155.8667
$ws.Range("M100").Value = -614.8667
# Row 122
$ws.Range("H122").Value = 7794.8076
$ws.Range("I122").Value = 7774
$ws.Range("J122").Value = 7841.625
$ws.Range("K122").Value = 23322
$ws.Range("L122").Value = 23524.875
$ws.Range("M122").Value = -20872
$ws.Range("N122").Value = -28424.875
# Row 126
$ws.Range("H126").Value = 44516588
$ws.Range("I126").Value = 46539990
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 139619970
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -139617500
$ws.Range("N126").Value = -10040
